$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "cliente" (sheet1.xml)
#   - new row 11: Id=1, Nome="Jailson Silva", Telefone=9991999222, Cpf="480.222.000-27"
#   - column C width 10 -> 11, new column D width 12
#   - selection moves to H6
# ---------------------------------------------------------------------------
$wsCliente = $wb.Worksheets.Item("cliente")
$wsCliente.Range("A11").Value = 1
$wsCliente.Range("B11").Value = "Jailson Silva"
$wsCliente.Range("C11").Value = 9991999222
$wsCliente.Range("D11").Value = "480.222.000-27"

$wsCliente.Columns.Item(3).ColumnWidth = 10.166666666666666
$wsCliente.Columns.Item(4).ColumnWidth = 11.166666666666666

$wsCliente.Range("H6").Select()

# ---------------------------------------------------------------------------
# Sheet "quarto" (sheet2.xml)
#   - B8 price 270 -> 270.5
#   - selection moves to B8
# ---------------------------------------------------------------------------
$wsQuarto = $wb.Worksheets.Item("quarto")
$wsQuarto.Range("B8").Value = 270.5

$wsQuarto.Range("B8").Select()

# ---------------------------------------------------------------------------
# Sheet "data" (sheet4.xml)
#   - A2:A8 change from text dates (shared strings) to real Excel date values
#   - column A width 10.42578125 -> 10.7109375 (best achievable approximation)
#   - selection moves to C2
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

$wsData.Range("A2").Value = (Get-Date -Year 2025 -Month 5 -Day 1).Date
$wsData.Range("A3").Value = (Get-Date -Year 2025 -Month 5 -Day 13).Date
$wsData.Range("A4").Value = (Get-Date -Year 2025 -Month 5 -Day 14).Date
$wsData.Range("A5").Value = (Get-Date -Year 2025 -Month 5 -Day 12).Date
$wsData.Range("A6").Value = (Get-Date -Year 2025 -Month 5 -Day 15).Date
$wsData.Range("A7").Value = (Get-Date -Year 2025 -Month 5 -Day 16).Date
$wsData.Range("A8").Value = (Get-Date -Year 2025 -Month 5 -Day 17).Date

$wsData.Columns.Item(1).ColumnWidth = 9.875

$wsData.Range("C2").Select()

# ---------------------------------------------------------------------------
# Sheet "reserva" (sheet5.xml)
#   - add explicit column widths for B:G
#   - selection moves to G10
# ---------------------------------------------------------------------------
$wsReserva = $wb.Worksheets.Item("reserva")

$wsReserva.Columns.Item(2).ColumnWidth = 13.877604166666666
$wsReserva.Columns.Item(3).ColumnWidth = 15.022135416666666
$wsReserva.Columns.Item(4).ColumnWidth = 25.592447916666668
$wsReserva.Columns.Item(5).ColumnWidth = 13.166666666666666
$wsReserva.Columns.Item(6).ColumnWidth = 19.451822916666668
$wsReserva.Columns.Item(7).ColumnWidth = 9.592447916666666

$wsReserva.Range("G10").Select()

# Restore "cliente" as the active tab (it was the active sheet before the
# edits and is not touched by the diff's tabSelected attribute), while each
# sheet keeps its own last-set selection from above.
$wsCliente.Activate()

Write-Host "edit complete"
